$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.932.03'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.817.55'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.35%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.01'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4681'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.37%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3703'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.54%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07385'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8727'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.05%  '

$ws.Range("E11").Value = '  -0.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.861.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.78%  '

$ws.Range("E13").Value = '  -0.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07076'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.516'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.90%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008734'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.960.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.339'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.047.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.902'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.56'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.230'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.90%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.329'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.26%  '

$ws.Range("E30").Value = '  -1.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08929'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7698'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.166'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.499'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.911'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.092'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01964'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.35%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05286'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.85%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.963'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.15%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.293'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.83%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5368'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.384'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.86%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1668'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.473'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4957'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.30%  '

$ws.Range("E47").Value = '  +1.66%  '

$ws.Range("E48").Value = '  +0.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9999'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06297'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.54%  '
